$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("Y1").Value = "pill"
$ws.Range("A7").NumberFormat = "@"
$ws.Range("A7").Value = "424"
$ws.Range("G7").Value = 44484
$ws.Range("K7").Value = "heroin; fentanyl"
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 1
$ws.Range("T7").Value = 0
$ws.Range("V7").Value = 0
$ws.Range("W7").Value = "powder"
$ws.Range("Y7").Value = ""
$ws.Range("Z7").Value = 1
$ws.Range("AP7").Value = ""
$ws.Range("AQ7").Value = "unknown"
$ws.Range("AR7").Value = ""
$ws.Range("AV7").Value = 44603
$ws.Range("AX7").Value = 0
$ws.Range("AZ7").Value = 7
$ws.Range("BA7").Value = 7
$ws.Range("BB7").Value = 1
$ws.Range("BC7").Value = 1
$ws.Range("BD7").Value = 1
$ws.Range("BE7").Value = 1
$ws.Range("BH7").Value = 0
$ws.Range("BK7").Value = 0
$ws.Range("BW7").Value = 1
$ws.Range("CB7").Value = 0
$ws.Range("CE7").Value = 1
$ws.Range("CF7").Value = 1
$ws.Range("CG7").Value = 0
$ws.Range("A8").NumberFormat = "@"
$ws.Range("A8").Value = "280"
$ws.Range("G8").Value = 44701
$ws.Range("K8").Value = "MDMA"
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = 0
$ws.Range("T8").Value = 1
$ws.Range("V8").Value = 1
$ws.Range("W8").Value = "pill"
$ws.Range("Y8").Value = 1
$ws.Range("Z8").Value = ""
$ws.Range("AP8").Value = "?"
$ws.Range("AQ8").Value = "not overdose related"
$ws.Range("AR8").Value = 0
$ws.Range("AV8").Value = 44721
$ws.Range("AX8").Value = 1
$ws.Range("AZ8").Value = 3
$ws.Range("BA8").Value = 1
$ws.Range("BB8").Value = 0
$ws.Range("BC8").Value = 0
$ws.Range("BD8").Value = 0
$ws.Range("BE8").Value = 0
$ws.Range("BH8").Value = 1
$ws.Range("BK8").Value = 1
$ws.Range("BW8").Value = 0
$ws.Range("CB8").Value = 1
$ws.Range("CE8").Value = 0
$ws.Range("CF8").Value = 0
$ws.Range("CG8").Value = 1
$ws.Range("U10").Value = "white"
$ws.Range("H12").Value = "spatula"
$ws.Range("I12").Value = "spatula"
$ws.Range("U12").Value = "white"
$ws.Range("A13").NumberFormat = "@"
$ws.Range("A13").Value = "253"
$ws.Range("G13").Value = 44707
$ws.Range("H13").Value = "swab"
$ws.Range("I13").Value = "swab"
$ws.Range("K13").Value = "crack"
$ws.Range("L13").Value = 0
$ws.Range("M13").Value = 0
$ws.Range("O13").Value = 1
$ws.Range("R13").Value = 1
$ws.Range("U13").Value = "white"
$ws.Range("AV13").Value = 44714
$ws.Range("AZ13").Value = 7
$ws.Range("BA13").Value = 1
$ws.Range("BB13").Value = 0
$ws.Range("BC13").Value = 0
$ws.Range("BD13").Value = 0
$ws.Range("BE13").Value = 0
$ws.Range("BI13").Value = 1
$ws.Range("BJ13").Value = 1
$ws.Range("BM13").Value = 1
$ws.Range("BW13").Value = 0
$ws.Range("CA13").Value = 1
$ws.Range("CC13").Value = 0
$ws.Range("CE13").Value = 0
$ws.Range("A14").NumberFormat = "@"
$ws.Range("A14").Value = "299"
$ws.Range("G14").Value = 44695
$ws.Range("H14").Value = "swab; spatula"
$ws.Range("I14").Value = "multiple methods"
$ws.Range("K14").Value = "heroin; fentanyl"
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 1
$ws.Range("O14").Value = 0
$ws.Range("R14").Value = 0
$ws.Range("U14").Value = "light gray"
$ws.Range("AV14").Value = 44707
$ws.Range("AZ14").Value = 9
$ws.Range("BA14").Value = 5
$ws.Range("BB14").Value = 1
$ws.Range("BC14").Value = 1
$ws.Range("BD14").Value = 1
$ws.Range("BE14").Value = 1
$ws.Range("BI14").Value = 0
$ws.Range("BJ14").Value = 0
$ws.Range("BM14").Value = 0
$ws.Range("BW14").Value = 1
$ws.Range("CA14").Value = 0
$ws.Range("CC14").Value = 1
$ws.Range("CE14").Value = 1
$ws.Range("U15").Value = "white"
$ws.Range("H18").Value = "pill"
$ws.Range("I18").Value = "pill"
$ws.Range("U18").Value = "blue"
$ws.Range("H20").Value = "spatula"
$ws.Range("I20").Value = "spatula"
$ws.Range("U20").Value = "white"
